$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated Cypher query text for the "FilesTab" row (cell B4): the "File Type"
# and "Breed" columns were dropped from the RETURN clause.
$query = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
WHERE samp.summarized_sample_type IN ["Primary Malignant Tumor Tissue"] 
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $query

# The user's last selection before saving moved from D4 to B4.
$ws.Range("B4").Select() | Out-Null
